# Weekly update: insert a new daily price record for
# Fruta / Terminal Hortofrutícola Agro Chillán - Mango, shifting the
# existing rows (old row 57 onward) down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 57; Excel shifts rows 57:81 down to 58:82 and
# copies formatting (incl. the date-style on column D) from the row above.
$ws.Rows(57).Insert()

# Populate the newly inserted row 57 with the new weekly record.
$ws.Cells.Item(57, 1).Value = 7
$ws.Cells.Item(57, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(57, 3).Value = "Ñuble"
$ws.Cells.Item(57, 4).Value = 44841
$ws.Cells.Item(57, 5).Value = 16
$ws.Cells.Item(57, 6).Value = "Fruta"
$ws.Cells.Item(57, 7).Value = 100108
$ws.Cells.Item(57, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(57, 9).Value = 100108002
$ws.Cells.Item(57, 10).Value = "Mango"
$ws.Cells.Item(57, 11).Value = "Sin especificar"
$ws.Cells.Item(57, 12).Value = "Primera"
$ws.Cells.Item(57, 13).Value = 80
$ws.Cells.Item(57, 14).Value = 7500
$ws.Cells.Item(57, 15).Value = 8000
$ws.Cells.Item(57, 16).Value = 7750
$ws.Cells.Item(57, 17).Value = "`$/bandeja 4 kilos"
$ws.Cells.Item(57, 18).Value = "Brasil"
$ws.Cells.Item(57, 19).Value = 1938
$ws.Cells.Item(57, 20).Value = 4
